# Add a new "truncateFrames" parameter row to the functional pipeline
# parameter sheet, just above the existing "scrubbing" section (new row 50,
# pushing the scrubbing.* / saveTimeSeries rows down by one).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 50 (shifts rows 50-56 down to 51-57).
$ws.Rows.Item(50).Insert()

# Populate the new row with the truncateFrames parameter definition.
$ws.Range("A50").Value = "reconstruction_functional_network.truncateFrames"
$ws.Range("D50").Value = "reconstruction_functional_network"
$ws.Range("E50").Value = "numeric"
$ws.Range("F50").Value = "scalar nonempty integer nonnegative"
$ws.Range("G50").Value = "standard"
$ws.Range("H50").Value = "Number of frames at beginning and end of timeseries to always remove"

# Update view state to match where the author ended up after the edit.
$ws.Range("F51").Select()
